# Refresh cryptos price/volume snapshot (scheduled GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text first so numeric-looking strings
# (e.g. "1.001", "15.00") are not auto-coerced into numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.071.06"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "1.649.72"
$ws.Range("E3").Value = "  -5.55%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "237.57"
$ws.Range("E5").Value = "  -4.73%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4796"
$ws.Range("E7").Value = "  -7.06%  "
$ws.Range("D8").Value = "0.2613"
$ws.Range("E8").Value = "  -5.33%  "
$ws.Range("D9").Value = "0.06004"
$ws.Range("E9").Value = "  -3.22%  "
$ws.Range("D10").Value = "0.07166"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "1.654.69"
$ws.Range("E11").Value = "  -5.17%  "
$ws.Range("D12").Value = "14.79"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").Value = "0.6232"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("D14").Value = "4.602"
$ws.Range("D15").Value = "73.38"
$ws.Range("E15").Value = "  -5.82%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "25.069.43"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").Value = "0.000006601"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "4.480"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").Value = "1.863.18"
$ws.Range("E22").Value = "  -5.22%  "
$ws.Range("D23").Value = "8.614"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "5.298"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "132.27"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "15.00"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").Value = "1.391"
$ws.Range("E27").Value = "  -7.22%  "
$ws.Range("D28").Value = "103.32"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").Value = "1.682"
$ws.Range("E29").Value = "  -5.90%  "
$ws.Range("D30").Value = "3.768"
$ws.Range("E30").Value = "  -5.09%  "
$ws.Range("D31").Value = "0.07911"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("D32").Value = "3.578"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").Value = "0.04596"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").Value = "2.598"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").Value = "0.9434"
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("D36").Value = "0.5773"
$ws.Range("E36").Value = "  -7.51%  "
$ws.Range("D37").Value = "2.620"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("D38").Value = "0.01552"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "0.8260"
$ws.Range("E40").Value = "  +8.99%  "
$ws.Range("D41").Value = "1.830"
$ws.Range("E41").Value = "  -5.58%  "
$ws.Range("D42").Value = "98.69"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Value = "0.3726"
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("D44").Value = "4.819"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").Value = "0.1143"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "6.109"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").Value = "0.05184"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "29.83"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").Value = "51.00"
$ws.Range("E49").Value = "  -8.15%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "0.3344"
$ws.Range("E51").Value = "  -2.98%  "

# Restore the default cell style so formatting matches the original file
# (values remain text; only the transient "@" number format is undone).
$dataRange.Style = "Normal"

